# feat: add 2022-Q1 data
#
# Before:
#   Sheet1 "2021-Q2"  (fund snapshot table)
#   Sheet2 "总计"      (date / holding-count / holding-value summary table)
#
# After:
#   Sheet1 "2021-Q2"   (unchanged)
#   Sheet2 "2022-Q1"   (new fund snapshot table, same shape as "2021-Q2")
#   Sheet3 "总计"      (summary table, with a new first row for 2022-Q1)

$wb = $excel.ActiveWorkbook

$wsTotalOld = $wb.Worksheets.Item(2)          # currently "总计", will become "2022-Q1"

# Grab a cell that already carries the bordered/bold/centered header style (s=2)
# so we can stamp it onto the newly created cells further down.
$headerStyleSrc = $wsTotalOld.Range("B1")

# --- create the brand new "总计" sheet right after the old one, while the old ---
# --- one still holds its original data/formatting.                           ---
$wsTotalNew = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTotalOld)

# Match the page margins used by the rest of the workbook (0.75/0.75/1/1/0.5/0.5 in).
$wsTotalNew.PageSetup.LeftMargin = 54
$wsTotalNew.PageSetup.RightMargin = 54
$wsTotalNew.PageSetup.TopMargin = 72
$wsTotalNew.PageSetup.BottomMargin = 72
$wsTotalNew.PageSetup.HeaderMargin = 36
$wsTotalNew.PageSetup.FooterMargin = 36

$headerStyleSrc.Copy() | Out-Null
$wsTotalNew.Range("B1:D1").PasteSpecial(-4122)
$wsTotalNew.Range("A2:A3").PasteSpecial(-4122)

$wsTotalNew.Range("B1").Value = "日期"
$wsTotalNew.Range("C1").Value = "持有数量(只)"
$wsTotalNew.Range("D1").Value = "持有市值(亿元)"

$wsTotalNew.Range("A2").Value = 0
$wsTotalNew.Range("B2").Value = "2022-Q1"
$wsTotalNew.Range("C2").Value = 1
$wsTotalNew.Range("D2").Value = 0

$wsTotalNew.Range("A3").Value = 1
$wsTotalNew.Range("B3").Value = "2021-Q2"
$wsTotalNew.Range("C3").Value = 1
$wsTotalNew.Range("D3").Value = 0

# --- turn the old "总计" sheet into the new "2022-Q1" fund snapshot sheet ---
$headerStyleSrc.Copy() | Out-Null
$wsTotalOld.Range("E1:H1").PasteSpecial(-4122)

$wsTotalOld.Range("B1").Value = "基金代码"
$wsTotalOld.Range("C1").Value = "基金名称"
$wsTotalOld.Range("D1").Value = "基金规模"
$wsTotalOld.Range("E1").Value = "股票总仓位"
$wsTotalOld.Range("F1").Value = "仓位占比"
$wsTotalOld.Range("G1").Value = "持有市值(亿元)"
$wsTotalOld.Range("H1").Value = "仓位排名"

$wsTotalOld.Range("A2").Value = 0

$wsTotalOld.Range("B2").NumberFormat = "@"
$wsTotalOld.Range("B2").Value = "000049"
$wsTotalOld.Range("B2").Style = "Normal"

$wsTotalOld.Range("C2").Value = "中银标普全球精选自然资源等权重指数(QDII)"

$wsTotalOld.Range("D2").NumberFormat = "@"
$wsTotalOld.Range("D2").Value = "0.27"
$wsTotalOld.Range("D2").Style = "Normal"

$wsTotalOld.Range("E2").NumberFormat = "@"
$wsTotalOld.Range("E2").Value = "89.72"
$wsTotalOld.Range("E2").Style = "Normal"

$wsTotalOld.Range("F2").NumberFormat = "@"
$wsTotalOld.Range("F2").Value = "1.13"
$wsTotalOld.Range("F2").Style = "Normal"

$wsTotalOld.Range("G2").NumberFormat = "@"
$wsTotalOld.Range("G2").Value = "0.0031"
$wsTotalOld.Range("G2").Style = "Normal"

$wsTotalOld.Range("H2").Value = 10

$wsTotalOld.Name = "2022-Q1"
$wsTotalNew.Name = "总计"

# Restore the originally active/selected tab ("2021-Q2").
$wb.Worksheets.Item(1).Activate()
